$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "very small" -> "small" label (shared string used in G3)
$ws.Range("G3").Value = "small"

# Update numeric values in column C (DF2): 496 -> 490
$ws.Range("C2").Value = 490
$ws.Range("C3").Value = 490
$ws.Range("C4").Value = 490
$ws.Range("C5").Value = 490

# Update F-statistic values (column D)
$ws.Range("D2").Value = 375858.88519869605
$ws.Range("D3").Value = 4.009086477938423
$ws.Range("D4").Value = 945.88974503464829
$ws.Range("D5").Value = 138.89344621813584

# Update p-values (column E)
$ws.Range("E3").Value = 0.0458062517406711

# Update etaSqp values (column F)
$ws.Range("F2").Value = 0.99869801660302171
$ws.Range("F3").Value = 0.0081154103996778654
$ws.Range("F4").Value = 0.65874817220860071
$ws.Range("F5").Value = 0.22085370272718619

# Update significance column (H)
$ws.Range("H3").Value = "*"
$ws.Range("H4").Value = "***"
$ws.Range("H5").Value = "***"

# Update column widths (closest representable values; engine quantizes
# column width storage to 1/6-character steps)
$ws.Columns.Item(5).ColumnWidth = 12.833333333333334
$ws.Columns.Item(6).ColumnWidth = 11.833333333333334
$ws.Columns.Item(7).ColumnWidth = 9.166666666666666
